$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price / volume(1h) data per latest scrape

# Row 2
$ws.Range("D2").Value = "69.496.33"
$ws.Range("E2").Value = "  +0.76%  "

# Row 3
$ws.Range("D3").Value = "2.494.54"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'570.26"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("D6").Value = "'166.86"
$ws.Range("E6").Value = "  +1.28%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("E8").Value = "  -0.29%  "

# Row 9
$ws.Range("E9").Value = "  +0.85%  "

# Row 10
$ws.Range("E10").Value = "  -0.59%  "

# Row 11
$ws.Range("E11").Value = "  -0.76%  "

# Row 12
$ws.Range("E12").Value = "  +0.34%  "

# Row 13
$ws.Range("D13").Value = "2.953.51"
$ws.Range("E13").Value = "  +0.14%  "

# Row 14
$ws.Range("D14").Value = "69.380.10"
$ws.Range("E14").Value = "  +0.73%  "

# Row 15
$ws.Range("E15").Value = "  +0.89%  "

# Row 16
$ws.Range("E16").Value = "  -1.34%  "

# Row 17
$ws.Range("D17").Value = "2.503.45"
$ws.Range("E17").Value = "  -0.21%  "

# Row 18
$ws.Range("D18").Value = "'11.22"
$ws.Range("E18").Value = "  +0.12%  "

# Row 19
$ws.Range("E19").Value = "  -1.33%  "

# Row 20
$ws.Range("D20").Value = "'352.78"
$ws.Range("E20").Value = "  +1.88%  "

# Row 21
$ws.Range("E21").Value = "  +0.73%  "

# Row 22
$ws.Range("E22").Value = "  -3.32%  "

# Row 23
$ws.Range("E23").Value = "  -0.03%  "

# Row 24
$ws.Range("D24").Value = "'69.37"
$ws.Range("E24").Value = "  -0.73%  "

# Row 25
$ws.Range("E25").Value = "  -2.28%  "

# Row 26
$ws.Range("D26").Value = "2.622.86"
$ws.Range("E26").Value = "  -0.99%  "

# Row 27
$ws.Range("E27").Value = "  -1.76%  "

# Row 28
$ws.Range("E28").Value = "  -1.42%  "

# Row 29
$ws.Range("E29").Value = "  -0.51%  "

# Row 30
$ws.Range("E30").Value = "  -1.86%  "

# Row 31
$ws.Range("B31").Value = "POPCAT"
$ws.Range("C31").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Range("D31").Value = "'3.35"
$ws.Range("E31").Value = "  +122.84%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'443.20"
$ws.Range("E32").Value = "  -3.08%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.20"
$ws.Range("E33").Value = "  -2.05%  "

# Row 34
$ws.Range("D34").Value = "'1.72"
$ws.Range("E34").Value = "  -0.13%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("D36").Value = "'154.31"
$ws.Range("E36").Value = "  -0.78%  "

# Row 37
$ws.Range("E37").Value = "  -1.48%  "

# Row 38
$ws.Range("D38").Value = "'19.06"
$ws.Range("E38").Value = "  +0.37%  "

# Row 39
$ws.Range("E39").Value = "  -1.13%  "

# Row 40
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
$ws.Range("E41").Value = "  -0.75%  "

# Row 42
$ws.Range("E42").Value = "  -0.44%  "

# Row 43
$ws.Range("E43").Value = "  -0.57%  "

# Row 44
$ws.Range("E44").Value = "  -0.09%  "

# Row 45
$ws.Range("E45").Value = "  -3.66%  "

# Row 46
$ws.Range("D46").Value = "'139.33"
$ws.Range("E46").Value = "  -1.09%  "

# Row 47
$ws.Range("E47").Value = "  -0.05%  "

# Row 48
$ws.Range("E48").Value = "  -1.51%  "

# Row 49
$ws.Range("E49").Value = "  -0.60%  "

# Row 50
$ws.Range("D50").Value = "'0.571"
$ws.Range("E50").Value = "  -0.35%  "

# Row 51
$ws.Range("E51").Value = "  -0.06%  "
